$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells that hold numeric-looking text stay as Text,
# matching the source data which uses "." as a thousands separator (not a valid number).
$priceTextCells = @("D2", "D3", "D5", "D6", "D9", "D11", "D12", "D14", "D15", "D16", "D17", "D19", "D22", "D23", "D24", "D28", "D30", "D31", "D32", "D33", "D35", "D38", "D39", "D43", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "39.355.74"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "2.193.81"
$ws.Range("E3").Value = "  -6.57%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "295.16"
$ws.Range("E5").Value = "  -4.54%  "
$ws.Range("D6").Value = "81.77"
$ws.Range("E7").Value = "  -3.95%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "0.468"
$ws.Range("E9").Value = "  -4.25%  "
$ws.Range("E10").Value = "  -6.83%  "
$ws.Range("D11").Value = "29.07"
$ws.Range("E11").Value = "  -4.77%  "
$ws.Range("D12").Value = "47.03"
$ws.Range("E12").Value = "  -11.02%  "
$ws.Range("E13").Value = "  -2.61%  "
$ws.Range("D14").Value = "2.532.84"
$ws.Range("E14").Value = "  -6.49%  "
$ws.Range("D15").Value = "6.23"
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("D16").Value = "13.90"
$ws.Range("E16").Value = "  -6.75%  "
$ws.Range("D17").Value = "2.195.86"
$ws.Range("E17").Value = "  -6.46%  "
$ws.Range("E18").Value = "  -6.42%  "
$ws.Range("D19").Value = "39.233.58"
$ws.Range("E19").Value = "  -2.41%  "
$ws.Range("E20").Value = "  -4.53%  "
$ws.Range("E21").Value = "  -6.87%  "
$ws.Range("D22").Value = "64.74"
$ws.Range("E22").Value = "  -4.82%  "
$ws.Range("D23").Value = "10.28"
$ws.Range("E23").Value = "  -5.03%  "
$ws.Range("D24").Value = "225.42"
$ws.Range("E24").Value = "  -4.68%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("E26").Value = "  -7.05%  "
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("D28").Value = "22.56"
$ws.Range("E28").Value = "  -4.64%  "
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").Value = "9.07"
$ws.Range("E30").Value = "  -2.43%  "
$ws.Range("D31").Value = "148.47"
$ws.Range("E31").Value = "  -2.25%  "
$ws.Range("D32").Value = "31.80"
$ws.Range("E32").Value = "  -9.33%  "
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  -7.50%  "
$ws.Range("D35").Value = "0.0693"
$ws.Range("E35").Value = "  -4.85%  "
$ws.Range("E36").Value = "  -4.78%  "
$ws.Range("E37").Value = "  -3.60%  "
$ws.Range("D38").Value = "15.32"
$ws.Range("E38").Value = "  -3.98%  "
$ws.Range("D39").Value = "0.0954"
$ws.Range("E39").Value = "  -5.50%  "
$ws.Range("E40").Value = "  -6.35%  "
$ws.Range("E41").Value = "  -5.01%  "
$ws.Range("E42").Value = "  -6.18%  "
$ws.Range("D43").Value = "1.900.64"
$ws.Range("E43").Value = "  -2.82%  "
$ws.Range("E44").Value = "  -9.15%  "
$ws.Range("E45").Value = "  -4.17%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "8.94"
$ws.Range("E46").Value = "  -4.67%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "16.01"
$ws.Range("E47").Value = "  -10.36%  "
$ws.Range("E48").Value = "  -4.55%  "
$ws.Range("D49").Value = "71.78"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "2.403.17"
$ws.Range("E50").Value = "  -6.23%  "
$ws.Range("D51").Value = "87.07"
$ws.Range("E51").Value = "  -6.93%  "
